$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-key the existing order row (row 2) with a new OrderID, and append two
#    more order rows (3 & 4) for the new orders created in ops.
# ---------------------------------------------------------------------------

# Row 2 - existing order, OrderID renamed BC18-001 -> 001BC
$ws.Range("A2").Value = 45436
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"
$ws.Range("B2").Value = "001BC"
$ws.Range("C2").Value = "SIPL5316"
$ws.Range("D2").Value = "SIPL5688"
$ws.Range("E2").Value = "SIPL0102"
$ws.Range("F2").Value = "SIPL0103"
$ws.Range("G2").Value = "BC Law Firm"
$ws.Range("H2").Value = "Title"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "Full Search"
$ws.Range("K2").Value = "FL"
$ws.Range("L2").Value = "Clay"
$ws.Range("M2").Value = "WIP"

# Row 3 - new order 002BC
$ws.Range("A3").Value = 45439
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = "002BC"
$ws.Range("C3").Value = "SIPL5316"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"
$ws.Range("G3").Value = "BC Law Firm"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Update Search"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "WIP"

# Row 4 - new order 003BC
$ws.Range("A4").Value = 45439
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = "003BC"
$ws.Range("C4").Value = "SIPL5316"
$ws.Range("D4").Value = "SIPL5688"
$ws.Range("E4").Value = "SIPL0102"
$ws.Range("F4").Value = "SIPL0103"
$ws.Range("G4").Value = "BC Law Firm"
$ws.Range("H4").Value = "Title"
$ws.Range("I4").Value = "Search & Typing"
$ws.Range("J4").Value = "Current Owner Search"
$ws.Range("K4").Value = "FL"
$ws.Range("L4").Value = "Clay"
$ws.Range("M4").Value = "WIP"

# ---------------------------------------------------------------------------
# 2. Match the font used across the data rows (B:M) to the explicit-black
#    Calibri font already used by column C/D in the template, and give the
#    date column (A) the same font with the new number format.
# ---------------------------------------------------------------------------
$ws.Range("A2:M4").Font.Color = 0

# ---------------------------------------------------------------------------
# 3. Re-style the header row: bold, explicit black font + gold fill.
# ---------------------------------------------------------------------------
$ws.Range("A1:M1").Font.Bold = $true
$ws.Range("A1:M1").Font.Color = 0
$ws.Range("A1:M1").Interior.Pattern = 1
$ws.Range("A1:M1").Interior.PatternColor = 0
$ws.Range("A1:M1").Interior.Color = 10086143

# ---------------------------------------------------------------------------
# 4. Column widths - let Excel re-fit the columns that now hold new/changed
#    text, mirroring what happens automatically after typing new data.
# ---------------------------------------------------------------------------
$ws.Columns("A:M").EntireColumn.AutoFit()

$ws.Range("G4").Select()
